$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "duration" column header
$ws.Range("C1").Value = "duration"

# Fill in duration values for each existing row
$ws.Range("C2").Value = 60
$ws.Range("C3").Value = 30
$ws.Range("C4").Value = 40
$ws.Range("C5").Value = 30

# Left-align the new numeric duration values (matches the new cell style)
$ws.Range("C2:C5").HorizontalAlignment = -4131

# Set the new column's width (closest achievable to the target 17.77734375)
$ws.Columns.Item(3).ColumnWidth = 17

# Update the active selection to match the edited workbook
$ws.Range("C5").Select()
